$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update "Förändrad" (column C) from 45184 to 45186 for every data row ---
# Data rows run from row 2 through row 199.
for ($r = 2; $r -le 199; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }
}

# --- 2) Add a friendly display-name second argument to every HYPERLINK() formula ---
# Columns: S=19 (artfynd), T=20 (kartor), U=21 (knärot), V=22 (klagomål),
#          W=23 (klagomålsmail), X=24 (tillsyn), Y=25 (tillsynsmail)
$hyperlinkCols = @(19, 20, 21, 22, 23, 24, 25)

for ($r = 2; $r -le 199; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq $null -or $label -eq "") {
        continue
    }
    foreach ($c in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -eq $null -or $f -eq "") {
            continue
        }
        if ($f.ToUpper().IndexOf("HYPERLINK(") -lt 0) {
            continue
        }
        # Already has two arguments? skip (nothing to do).
        if ($f.IndexOf(",") -ge 0) {
            continue
        }
        if ($f.Substring($f.Length - 1, 1) -eq ")") {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
